# "Pallets added, and CSV referencing"
#
# Adds a new "Pallets" block (Size / Approx Qty) next to the existing
# Blanks / Labor Rates summary on the "Totals" sheet, and switches that
# whole summary table's header formatting from Wrap-Text over to
# Shrink-to-Fit (with a tighter row 3) so everything lines up again now
# that a third block has been added.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Totals")

function Set-RangeStyle {
    param(
        [string]$addr,
        [string]$fontName = $null,
        $fontSize = $null,
        $bold = $null,
        $hAlign = $null,
        $vAlign = $null,
        $wrap = $null,
        $shrink = $null,
        [string]$numberFormat = $null
    )
    $rng = $ws.Range($addr)
    if ($fontName) { $rng.Font.Name = $fontName }
    if ($fontSize) { $rng.Font.Size = $fontSize }
    if ($null -ne $bold) { $rng.Font.Bold = $bold }
    if ($null -ne $hAlign) { $rng.HorizontalAlignment = $hAlign }
    if ($null -ne $vAlign) { $rng.VerticalAlignment = $vAlign }
    if ($null -ne $wrap) { $rng.WrapText = $wrap }
    if ($null -ne $shrink) { $rng.ShrinkToFit = $shrink }
    if ($numberFormat) { $rng.NumberFormat = $numberFormat }
}

$xlCenter = -4108

# ---------------------------------------------------------------------
# New values for the Pallets block
# ---------------------------------------------------------------------
$ws.Range("G1").Value = "Pallets"
$ws.Range("G3").Value = "Size"
$ws.Range("H3").Value = "Approx Qty"
$ws.Range("G1:H1").Merge()

# ---------------------------------------------------------------------
# Row 1 header groups (Blanks / Labor Rates / Pallets): all three now
# share one consistent look - bold 16pt, centered, shrink-to-fit.
# ---------------------------------------------------------------------
foreach ($addr in @("A1:B1", "D1:E1", "G1:H1")) {
    Set-RangeStyle -addr $addr -hAlign $xlCenter -vAlign $xlCenter -wrap $false -shrink $true -numberFormat "@"
}

# ---------------------------------------------------------------------
# Row 2 (blank spacer row under the header groups): shrink-to-fit too
# ---------------------------------------------------------------------
foreach ($addr in @("A2:B2", "E2")) {
    Set-RangeStyle -addr $addr -hAlign $xlCenter -vAlign $xlCenter -wrap $false -shrink $true
}

# ---------------------------------------------------------------------
# Row 3 (column headers): style the two new header cells like the rest
# of the row, switch the whole row from Wrap to Shrink, and shrink the
# row back down now it no longer needs the extra height for wrapping.
# ---------------------------------------------------------------------
Set-RangeStyle -addr "G3:H3" -fontName "Calibri" -fontSize 12 -bold $true -hAlign $xlCenter -vAlign $xlCenter -shrink $true -numberFormat "@"
foreach ($addr in @("A3", "B3", "D3", "E3")) {
    Set-RangeStyle -addr $addr -hAlign $xlCenter -vAlign $xlCenter -wrap $false -shrink $true
}
$ws.Rows.Item(3).RowHeight = 16.5

# Narrow spacer cells either side of the new Pallets block (same 12pt
# font as the rest of the row, just shrink-to-fit with no forced align)
foreach ($addr in @("C3", "F3")) {
    Set-RangeStyle -addr $addr -fontName "Calibri" -fontSize 12 -shrink $true
}

# ---------------------------------------------------------------------
# Rows 4-6 (Heian/Weeke, Vector, Pack totals): shrink-to-fit
# ---------------------------------------------------------------------
Set-RangeStyle -addr "D4:E6" -hAlign $xlCenter -vAlign $xlCenter -shrink $true

# ---------------------------------------------------------------------
# Column layout: re-proportion the existing columns, add narrow spacer
# columns C & F, and size the two new Pallets columns G & H.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 21        # A - Blanks
$ws.Columns.Item(2).ColumnWidth = 9.43      # B
$ws.Columns.Item(3).ColumnWidth = 1.57      # C - spacer
$ws.Columns.Item(5).ColumnWidth = 11.57     # E - Labor Rates value
$ws.Columns.Item(6).ColumnWidth = 1.43      # F - spacer
$ws.Columns.Item(7).ColumnWidth = 20.43     # G - Pallets
$ws.Columns.Item(8).ColumnWidth = 9.14      # H - Approx Qty

# ---------------------------------------------------------------------
# Selection matches where the editor left off after adding the block
# ---------------------------------------------------------------------
$ws.Range("G4").Select()
